$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.106.15'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '1.790.55'
$ws.Range('E3').Value = '  -1.06%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.17'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.552'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '32.51'
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0714'
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').Value = '2.047.25'
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.27'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').Value = '1.796.54'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.627'
$ws.Range('E15').Value = '  -3.25%  '
$ws.Range('D16').Value = '34.103.97'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('E17').Value = '  -3.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.10'
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.81'
$ws.Range('E19').Value = '  -3.98%  '
$ws.Range('E20').Value = '  -2.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.77'
$ws.Range('E22').Value = '  -1.01%  '
$ws.Range('E23').Value = '  -3.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.12'
$ws.Range('E24').Value = '  -1.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.27'
$ws.Range('E25').Value = '  -1.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.35'
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.05'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('E28').Value = '  -2.35%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  -3.39%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  -3.54%  '
$ws.Range('E33').Value = '  -4.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('E34').Value = '  -4.93%  '
$ws.Range('D35').Value = '1.385.62'
$ws.Range('E35').Value = '  -4.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.650'
$ws.Range('E36').Value = '  +0.91%  '
$ws.Range('E37').Value = '  -1.43%  '
$ws.Range('E38').Value = '  -3.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '80.01'
$ws.Range('E39').Value = '  -6.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.35'
$ws.Range('E40').Value = '  +0.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.918'
$ws.Range('E41').Value = '  -5.03%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.18'
$ws.Range('E42').Value = '  +1.39%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.69'
$ws.Range('E43').Value = '  -3.74%  '
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0499'
$ws.Range('E45').Value = '  +1.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '108.02'
$ws.Range('E46').Value = '  +1.91%  '
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('D48').Value = '1.947.56'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '12.04'
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0131'
$ws.Range('E50').Value = '  +4.43%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.07%  '
